# Weekly update: a new price record is inserted as row 17 (pushing the
# existing rows 17-47 down to 18-48), extending the used range to A1:R48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17; Excel shifts rows 17-47 down to 18-48
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with this week's record.
$ws.Cells.Item(17, 1).Value = 9
$ws.Cells.Item(17, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(17, 3).Value = "Metropolitana"
$ws.Cells.Item(17, 4).Value = 44580
$ws.Cells.Item(17, 5).Value = 13
$ws.Cells.Item(17, 6).Value = 100112029
$ws.Cells.Item(17, 7).Value = "Orégano"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 16
$ws.Cells.Item(17, 11).Value = 9000
$ws.Cells.Item(17, 12).Value = 10000
$ws.Cells.Item(17, 13).Value = 9500
$ws.Cells.Item(17, 14).Value = "`$/docena de atados"
$ws.Cells.Item(17, 15).Value = "Región Metropolitana"
$ws.Cells.Item(17, 16).Value = 3167
$ws.Cells.Item(17, 17).Value = 3
$ws.Cells.Item(17, 18).Value = "Hortaliza"
